$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 9, shifting the existing rows 9-41 down to rows 11-43.
$ws.Rows("9:10").Insert()

# New row 9: new observation for patrol 2 - "Between Dalupiri & Calayan Islands"
$ws.Cells.Item(9,1).Value = 2
$ws.Cells.Item(9,2).Value = 6
$ws.Cells.Item(9,3).Value = 323
$ws.Cells.Item(9,4).Value = -9
$ws.Cells.Item(9,5).Value = 16364
$ws.Cells.Item(9,6).Value = 19.216999999999999
$ws.Cells.Item(9,7).Value = 121.333
$ws.Cells.Item(9,8).Value = "Between Dalupiri & Calayan Islands"

# New row 10: new observation for patrol 2 - "North of Pratas Island"
$ws.Cells.Item(10,1).Value = 2
$ws.Cells.Item(10,2).Value = 7
$ws.Cells.Item(10,3).Value = 400
$ws.Cells.Item(10,4).Value = -9
$ws.Cells.Item(10,5).Value = 16366
$ws.Cells.Item(10,6).Value = 21.012
$ws.Cells.Item(10,7).Value = 116.751
$ws.Cells.Item(10,8).Value = "North of Pratas Island"

# The remaining patrol-2 observations (now at rows 11-14) were pushed down by the two
# newly inserted observations, so their sequence "number" column needs to be bumped by 1.
$ws.Cells.Item(11,2).Value = 7
$ws.Cells.Item(12,2).Value = 8
$ws.Cells.Item(13,2).Value = 9
$ws.Cells.Item(14,2).Value = 10

# Update the view: scroll so row 4 is at the top and select C11, matching the saved view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C11").Select()
